$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Battery Standby / Alarm Load test-data updates (NGC-488/T399 OR TC-151, 5/24/40V + AC calc rows)

# B4: test reference id (was blank)
$ws.Range("B4").Value = "NGC-488/T399 OR TC-151"

# C8: CPU type value cleared out (was "CPU 801") - typed as a literal blank with quote-prefix
$ws.Range("C8").Value = "'"

# F8: 40V standby current now recorded as text "0.300" (quote-prefixed number)
$ws.Range("F8").Value = "'0.300"

# K7/L7: Heat sensor alarm/standby current figures
$ws.Range("K7").Value = "'0.00"
$ws.Range("L7").Value = 0.005
# writing a plain number resets the quote-prefix formatting Excel had kept on this
# cell, so restore the original look (same format as the untouched K6/L6 cells)
# by copying formats only from L6 back onto L7.
$ws.Range("L6").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null

# K8/L8: Ion sensor alarm/standby current figures
$ws.Range("K8").Value = "'0.00"
$ws.Range("L8").Value = 0.005
$ws.Range("L6").Copy() | Out-Null
$ws.Range("L8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Match the author's final cursor position/selection on the sheet
$ws.Range("J8").Select() | Out-Null
